# Recomputed TPM-based NATMI ligand-receptor statistics for Fn1-Itgav.
# Updates average/total expression, specificity, and edge-weight columns (G:J, M:P, Q:T)
# on the "LR-pairs_lrc2p" sheet to reflect the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster = ECs, Target cluster = ECs
$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 336.7976575282865
$ws.Range("R2").Value = 3031.178917754578
$ws.Range("S2").Value = 0.001297415516188202
$ws.Range("T2").Value = 0.001297415516188201

# Row 3: Sending cluster = ECs, Target cluster = FAPs
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 1808.120147194053
$ws.Range("R3").Value = 16273.08132474648
$ws.Range("S3").Value = 0.006965259649720209
$ws.Range("T3").Value = 0.006965259649720207

# Row 4: Sending cluster = ECs, Target cluster = MuSCs
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 959.2580266172837
$ws.Range("R4").Value = 8633.322239555553
$ws.Range("S4").Value = 0.003695263966189588
$ws.Range("T4").Value = 0.003695263966189588

# Row 5: Sending cluster = ECs, Target cluster = Resolving-Mac
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 1957.150239956869
$ws.Range("R5").Value = 17614.35215961182
$ws.Range("S5").Value = 0.007539354957117662
$ws.Range("T5").Value = 0.007539354957117661

# Row 6: Sending cluster = FAPs, Target cluster = ECs
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 15919.61297430572
$ws.Range("R6").Value = 143276.5167687514
$ws.Range("S6").Value = 0.06132570231086169
$ws.Range("T6").Value = 0.06132570231086169

# Row 7: Sending cluster = FAPs, Target cluster = FAPs
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("Q7").Value = 85465.47848824183
$ws.Range("R7").Value = 769189.3063941764
$ws.Range("S7").Value = 0.3292310246539678
$ws.Range("T7").Value = 0.3292310246539678

# Row 8: Sending cluster = FAPs, Target cluster = MuSCs
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 45341.81335557788
$ws.Range("R8").Value = 408076.3202002009
$ws.Range("S8").Value = 0.1746662153512617
$ws.Range("T8").Value = 0.1746662153512617

# Row 9: Sending cluster = FAPs, Target cluster = Resolving-Mac
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 92509.77153861629
$ws.Range("R9").Value = 832587.9438475465
$ws.Range("S9").Value = 0.3563671252171522
$ws.Range("T9").Value = 0.3563671252171522

# Row 10: Sending cluster = MuSCs, Target cluster = ECs
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 882.4641677360071
$ws.Range("R10").Value = 7942.177509624064
$ws.Range("S10").Value = 0.00339943784675713
$ws.Range("T10").Value = 0.00339943784675713

# Row 11: Sending cluster = MuSCs, Target cluster = FAPs
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 4737.566325639595
$ws.Range("R11").Value = 42638.09693075635
$ws.Range("S11").Value = 0.01825010335571976
$ws.Range("T11").Value = 0.01825010335571976

# Row 12: Sending cluster = MuSCs, Target cluster = MuSCs
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 2513.410699811088
$ws.Range("R12").Value = 22620.69629829979
$ws.Range("S12").Value = 0.009682187413119037
$ws.Range("T12").Value = 0.009682187413119039

# Row 13: Sending cluster = MuSCs, Target cluster = Resolving-Mac
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 5128.049198183066
$ws.Range("R13").Value = 46152.4427836476
$ws.Range("S13").Value = 0.01975432562781526
$ws.Range("T13").Value = 0.01975432562781526

# Row 14: Sending cluster = Resolving-Mac, Target cluster = ECs
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 135.1969760912086
$ws.Range("R14").Value = 1216.772784820877
$ws.Range("S14").Value = 0.0005208072283214366
$ws.Range("T14").Value = 0.0005208072283214366

# Row 15: Sending cluster = Resolving-Mac, Target cluster = FAPs
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 725.8137663552371
$ws.Range("R15").Value = 6532.323897197134
$ws.Range("S15").Value = 0.00279598750550453
$ws.Range("T15").Value = 0.00279598750550453

# Row 16: Sending cluster = Resolving-Mac, Target cluster = MuSCs
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 385.0643898227964
$ws.Range("R16").Value = 3465.579508405168
$ws.Range("S16").Value = 0.001483349135365288
$ws.Range("T16").Value = 0.001483349135365289

# Row 17: Sending cluster = Resolving-Mac, Target cluster = Resolving-Mac
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 785.6372759247263
$ws.Range("R17").Value = 7070.735483322536
$ws.Range("S17").Value = 0.003026440264938493
$ws.Range("T17").Value = 0.003026440264938493
